# Rename "strategy_id-5008" to "strategy_id-5007", then insert a new sheet
# "strategy_id-5009" right after it, as a copy of its contents (both sheets
# hold the single-variable LHS template for frac_gnrl_eating_red_meat, which
# does not vary/sample across trajectories -- matches the commit's "Ls still
# vary for all Ls" fix and the PFLO:ALL_NO_STOPPING_DEFORESTATION_PLUR
# template rebuild).

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("strategy_id-5008")
$src.Name = "strategy_id-5007"

$src.Copy([System.Reflection.Missing]::Value, $src)

$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "strategy_id-5009"
